# This script applies a scheduled Market Board price refresh to the
# Chocobo Leve profit-tracking workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# For each affected leve row, columns H-L hold refreshed Universalis
# average-price data, and columns M/N hold the recomputed NQ/HQ leve
# profit (a column is omitted entirely whenever that profit figure is
# not applicable, e.g. exactly break-even/zero).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 28725.6
$ws.Range("J3").Value = 28725.6
$ws.Range("L3").Value = 28725.6
$ws.Range("N3").Value = -28953.6

$ws.Range("H6").Value = 760.625
$ws.Range("I6").Value = 120.76923
$ws.Range("K6").Value = 362.30769
$ws.Range("M6").Value = -250.30769

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H39").Value = 502.44446
$ws.Range("I39").Value = 224.4
$ws.Range("K39").Value = 673.2
$ws.Range("M39").Value = -377.2

$ws.Range("H102").Value = 28725.6
$ws.Range("J102").Value = 28725.6
$ws.Range("L102").Value = 28725.6
$ws.Range("N102").Value = -35215.6

$ws.Range("H129").Value = 1131.5333
$ws.Range("I129").Value = 316.66666
$ws.Range("J129").Value = 1222.0741
$ws.Range("K129").Value = 949.9999799999999
$ws.Range("L129").Value = 3666.2223
$ws.Range("M129").Value = 4050.00002
$ws.Range("N129").Value = -13666.2223

$ws.Range("H137").Value = 7332.5
$ws.Range("I137").Value = 8590.6
$ws.Range("K137").Value = 25771.8
$ws.Range("M137").Value = -23221.8

$ws.Range("H138").Value = 2406.61
$ws.Range("I138").Value = 718.619
$ws.Range("J138").Value = 2855.3164
$ws.Range("K138").Value = 2155.857
$ws.Range("L138").Value = 8565.949200000001
$ws.Range("M138").Value = 2984.143
$ws.Range("N138").Value = -18845.9492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5822.9683
$ws.Range("I32").Value = 3852.2888
$ws.Range("K32").Value = 3852.2888
$ws.Range("M32").Value = -3565.2888

$ws.Range("H34").Value = 21500
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2729

$ws.Range("H74").Value = 2873.2917
$ws.Range("I74").Value = 2835.525
$ws.Range("J74").Value = 3062.125
$ws.Range("K74").Value = 2835.525
$ws.Range("L74").Value = 3062.125
$ws.Range("M74").Value = -1961.525
$ws.Range("N74").Value = -4810.125

$ws.Range("H77").Value = 2873.2917
$ws.Range("I77").Value = 2835.525
$ws.Range("J77").Value = 3062.125
$ws.Range("K77").Value = 14177.625
$ws.Range("L77").Value = 15310.625
$ws.Range("M77").Value = -9809.625
$ws.Range("N77").Value = -24046.625

$ws.Range("H102").Value = 2414.2856
$ws.Range("I102").Value = 2350
$ws.Range("K102").Value = 2350
$ws.Range("M102").Value = -728

$ws.Range("H110").Value = 851.61536
$ws.Range("I110").Value = 782.1
$ws.Range("J110").Value = 1083.3334
$ws.Range("K110").Value = 782.1
$ws.Range("L110").Value = 1083.3334
$ws.Range("M110").Value = 1262.9
$ws.Range("N110").Value = -5173.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10206725
$ws.Range("I31").Value = 1399.9678
$ws.Range("J31").Value = 27782562
$ws.Range("K31").Value = 1399.9678
$ws.Range("L31").Value = 27782562
$ws.Range("M31").Value = -1104.9678
$ws.Range("N31").Value = -27783152

$ws.Range("H34").Value = 10206725
$ws.Range("I34").Value = 1399.9678
$ws.Range("J34").Value = 27782562
$ws.Range("K34").Value = 1399.9678
$ws.Range("L34").Value = 27782562
$ws.Range("M34").Value = -1197.9678
$ws.Range("N34").Value = -27782966

$ws.Range("H132").Value = 2623.9473
$ws.Range("I132").Value = 1087.238
$ws.Range("J132").Value = 4522.2354
$ws.Range("K132").Value = 3261.714
$ws.Range("L132").Value = 13566.7062
$ws.Range("M132").Value = -731.7139999999999
$ws.Range("N132").Value = -18626.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 705.875
$ws.Range("I113").Value = 620.54285
$ws.Range("J113").Value = 848.0952
$ws.Range("K113").Value = 1861.62855
$ws.Range("L113").Value = 2544.2856
$ws.Range("M113").Value = 308.3714499999999
$ws.Range("N113").Value = -6884.2856

$ws.Range("H121").Value = 1639.4559
$ws.Range("J121").Value = 1749.9048
$ws.Range("L121").Value = 5249.7144
$ws.Range("N121").Value = -7869.7144

$ws.Range("H139").Value = 1724.8125
$ws.Range("I139").Value = 1176.6923
$ws.Range("J139").Value = 4100
$ws.Range("K139").Value = 3530.0769
$ws.Range("L139").Value = 12300
$ws.Range("M139").Value = 1609.9231
$ws.Range("N139").Value = -22580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4219243.5
$ws.Range("I11").Value = 6545482
$ws.Range("J11").Value = 2086858.6
$ws.Range("K11").Value = 6545482
$ws.Range("L11").Value = 2086858.6
$ws.Range("M11").Value = -6545343
$ws.Range("N11").Value = -2087136.6

$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15346
$ws.Range("M21").ClearContents()

$ws.Range("H24").Value = 8999.5
$ws.Range("J24").Value = 8999.5
$ws.Range("L24").Value = 8999.5
$ws.Range("N24").Value = -9345.5

$ws.Range("H30").Value = 15000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 15000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15210
$ws.Range("M30").ClearContents()

$ws.Range("H132").Value = 4280.4
$ws.Range("I132").Value = 2744.5715
$ws.Range("J132").Value = 5107.385
$ws.Range("K132").Value = 8233.7145
$ws.Range("L132").Value = 15322.155
$ws.Range("M132").Value = -5703.7145
$ws.Range("N132").Value = -20382.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4628.85
$ws.Range("I7").Value = 2558.3
$ws.Range("J7").Value = 6699.4
$ws.Range("K7").Value = 2558.3
$ws.Range("L7").Value = 6699.4
$ws.Range("M7").Value = -2446.3
$ws.Range("N7").Value = -6923.4

$ws.Range("H100").Value = 2184.3333
$ws.Range("I100").Value = 2021.2
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2021.2
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1480.2
$ws.Range("N100").Value = -4082

$ws.Range("H126").Value = 4628.85
$ws.Range("I126").Value = 2558.3
$ws.Range("J126").Value = 6699.4
$ws.Range("K126").Value = 7674.900000000001
$ws.Range("L126").Value = 20098.2
$ws.Range("M126").Value = -5204.900000000001
$ws.Range("N126").Value = -25038.2

$ws.Range("H132").Value = 2955.5232
$ws.Range("I132").Value = 1635.1842
$ws.Range("K132").Value = 4905.5526
$ws.Range("M132").Value = -2375.5526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 9999.5
$ws.Range("J12").Value = 9999.5
$ws.Range("L12").Value = 9999.5
$ws.Range("N12").Value = -10283.5

$ws.Range("H14").Value = 1700
$ws.Range("I14").Value = 800
$ws.Range("J14").Value = 3500
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 3500
$ws.Range("M14").Value = -632
$ws.Range("N14").Value = -3836

$ws.Range("H19").Value = 7799.5
$ws.Range("J19").Value = 7799.5
$ws.Range("L19").Value = 7799.5
$ws.Range("N19").Value = -8147.5

$ws.Range("H39").Value = 18333
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4587

$ws.Range("H42").Value = 45024.5
$ws.Range("J42").Value = 45024.5
$ws.Range("L42").Value = 45024.5
$ws.Range("N42").Value = -45780.5

$ws.Range("H43").Value = 18405.8
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 28009.666
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 28009.666
$ws.Range("M43").Value = -3851
$ws.Range("N43").Value = -28307.666

$ws.Range("H74").Value = 5874.5
$ws.Range("J74").Value = 5874.5
$ws.Range("L74").Value = 5874.5
$ws.Range("N74").Value = -7746.5

$ws.Range("H77").Value = 5874.5
$ws.Range("J77").Value = 5874.5
$ws.Range("L77").Value = 17623.5
$ws.Range("N77").Value = -26983.5

$ws.Range("H127").Value = 39886
$ws.Range("J127").Value = 39886
$ws.Range("L127").Value = 39886
$ws.Range("N127").Value = -49806

$ws.Range("H132").Value = 7938050.5
$ws.Range("I132").Value = 803.4400000000001
$ws.Range("J132").Value = 19610472
$ws.Range("K132").Value = 2410.32
$ws.Range("L132").Value = 58831416
$ws.Range("M132").Value = 119.6799999999998
$ws.Range("N132").Value = -58836476
